$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 499.66666
$ws.Range("I4").Value = 499.66666
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 499.66666
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -385.66666
$ws.Range("N4").ClearContents()
$ws.Range("H12").Value = 114
$ws.Range("I12").Value = 120
$ws.Range("J12").Value = 102
$ws.Range("K12").Value = 120
$ws.Range("L12").Value = 102
$ws.Range("M12").Value = 50
$ws.Range("N12").Value = -442
$ws.Range("H17").Value = 946
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 946
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 2838
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -3174
$ws.Range("H28").Value = 3491.0908
$ws.Range("I28").Value = 3449.5
$ws.Range("J28").Value = 3541
$ws.Range("K28").Value = 3449.5
$ws.Range("L28").Value = 3541
$ws.Range("M28").Value = -2964.5
$ws.Range("N28").Value = -4511
$ws.Range("H34").Value = 4998.091
$ws.Range("I34").Value = 4998.091
$ws.Range("K34").Value = 4998.091
$ws.Range("M34").Value = -4795.091
$ws.Range("H36").Value = 4998.091
$ws.Range("I36").Value = 4998.091
$ws.Range("K36").Value = 4998.091
$ws.Range("M36").Value = -4283.091
$ws.Range("H40").Value = 3816.4285
$ws.Range("I40").Value = 3381.3333
$ws.Range("K40").Value = 3381.3333
$ws.Range("M40").Value = -3206.3333
$ws.Range("H51").Value = 3000
$ws.Range("I51").Value = 3000
$ws.Range("K51").Value = 3000
$ws.Range("M51").Value = -2516
$ws.Range("H53").Value = 1187.8
$ws.Range("I53").Value = 1187.8
$ws.Range("K53").Value = 1187.8
$ws.Range("M53").Value = -550.8
$ws.Range("H62").Value = 2970
$ws.Range("I62").Value = 2837.875
$ws.Range("J62").Value = 3498.5
$ws.Range("K62").Value = 2837.875
$ws.Range("L62").Value = 3498.5
$ws.Range("M62").Value = -2213.875
$ws.Range("N62").Value = -4746.5
$ws.Range("H65").Value = 2970
$ws.Range("I65").Value = 2837.875
$ws.Range("J65").Value = 3498.5
$ws.Range("K65").Value = 14189.375
$ws.Range("L65").Value = 17492.5
$ws.Range("M65").Value = -11069.375
$ws.Range("N65").Value = -23732.5
$ws.Range("H76").Value = 4595.6
$ws.Range("I76").Value = 5328.3335
$ws.Range("K76").Value = 5328.3335
$ws.Range("M76").Value = -5013.3335
$ws.Range("H79").Value = 4595.6
$ws.Range("I79").Value = 5328.3335
$ws.Range("K79").Value = 5328.3335
$ws.Range("M79").Value = -4236.3335
$ws.Range("H86").Value = 999.3333
$ws.Range("I86").Value = 999.3333
$ws.Range("K86").Value = 999.3333
$ws.Range("M86").Value = 123.6667
$ws.Range("H88").Value = 1797.625
$ws.Range("J88").Value = 1711.5714
$ws.Range("L88").Value = 1711.5714
$ws.Range("N88").Value = -2523.5714
$ws.Range("H89").Value = 999.3333
$ws.Range("I89").Value = 999.3333
$ws.Range("K89").Value = 4996.6665
$ws.Range("M89").Value = 619.3334999999997
$ws.Range("H91").Value = 1797.625
$ws.Range("J91").Value = 1711.5714
$ws.Range("L91").Value = 1711.5714
$ws.Range("N91").Value = -4519.5714
$ws.Range("H92").Value = 1341.8125
$ws.Range("I92").Value = 639.3333
$ws.Range("J92").Value = 3449.25
$ws.Range("K92").Value = 639.3333
$ws.Range("L92").Value = 3449.25
$ws.Range("M92").Value = 608.6667
$ws.Range("N92").Value = -5945.25
$ws.Range("H98").Value = 1764.091
$ws.Range("I98").Value = 1303.25
$ws.Range("J98").Value = 2993
$ws.Range("K98").Value = 1303.25
$ws.Range("L98").Value = 2993
$ws.Range("M98").Value = 194.75
$ws.Range("N98").Value = -5989
$ws.Range("H106").Value = 2202.1428
$ws.Range("I106").Value = 2137.25
$ws.Range("K106").Value = 2137.25
$ws.Range("M106").Value = -1506.25
$ws.Range("H113").Value = 11499.286
$ws.Range("I113").Value = 10299
$ws.Range("K113").Value = 10299
$ws.Range("M113").Value = -7045
$ws.Range("H122").Value = 1764.091
$ws.Range("I122").Value = 1303.25
$ws.Range("J122").Value = 2993
$ws.Range("K122").Value = 3909.75
$ws.Range("L122").Value = 8979
$ws.Range("M122").Value = -1459.75
$ws.Range("N122").Value = -13879
$ws.Range("H132").Value = 209437.25
$ws.Range("I132").Value = 1044.2
$ws.Range("J132").Value = 3335333
$ws.Range("K132").Value = 3132.6
$ws.Range("L132").Value = 10005999
$ws.Range("M132").Value = -602.6000000000004
$ws.Range("N132").Value = -10011059
$ws.Range("H135").Value = 397.92856
$ws.Range("I135").Value = 397.92856
$ws.Range("K135").Value = 3581.35704
$ws.Range("M135").Value = -1046.35704
$ws.Range("H137").Value = 3235.7144
$ws.Range("I137").Value = 2719.4443
$ws.Range("J137").Value = 3622.9167
$ws.Range("K137").Value = 8158.3329
$ws.Range("L137").Value = 10868.7501
$ws.Range("M137").Value = -5608.3329
$ws.Range("N137").Value = -15968.7501
$ws.Range("H138").Value = 3714.5747
$ws.Range("I138").Value = 3376.1628
$ws.Range("J138").Value = 4045.2954
$ws.Range("K138").Value = 10128.4884
$ws.Range("L138").Value = 12135.8862
$ws.Range("M138").Value = -4988.4884
$ws.Range("N138").Value = -22415.8862

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 852.3333
$ws.Range("J26").Value = 500
$ws.Range("L26").Value = 500
$ws.Range("N26").Value = -1160
$ws.Range("H32").Value = 2330.7231
$ws.Range("I32").Value = 2055.5398
$ws.Range("K32").Value = 2055.5398
$ws.Range("M32").Value = -1768.5398
$ws.Range("H43").Value = 60404.332
$ws.Range("J43").Value = 30776.166
$ws.Range("L43").Value = 30776.166
$ws.Range("N43").Value = -31402.166
$ws.Range("H45").Value = 3588.4167
$ws.Range("I45").Value = 3010.8
$ws.Range("K45").Value = 3010.8
$ws.Range("M45").Value = -2633.8
$ws.Range("H61").Value = 2291.2903
$ws.Range("I61").Value = 2052.6155
$ws.Range("K61").Value = 2052.6155
$ws.Range("M61").Value = -1840.6155
$ws.Range("H74").Value = 1390.3846
$ws.Range("I74").Value = 1372.2174
$ws.Range("K74").Value = 1372.2174
$ws.Range("M74").Value = -498.2174
$ws.Range("H77").Value = 1390.3846
$ws.Range("I77").Value = 1372.2174
$ws.Range("K77").Value = 6861.087
$ws.Range("M77").Value = -2493.087
$ws.Range("H97").Value = 1228.9231
$ws.Range("I97").Value = 831.3333
$ws.Range("J97").Value = 6000
$ws.Range("K97").Value = 831.3333
$ws.Range("L97").Value = 6000
$ws.Range("M97").Value = -335.3333
$ws.Range("N97").Value = -6992
$ws.Range("H102").Value = 2069.875
$ws.Range("I102").Value = 1019.6667
$ws.Range("K102").Value = 1019.6667
$ws.Range("M102").Value = 602.3333
$ws.Range("H110").Value = 2297
$ws.Range("I110").Value = 2282.2856
$ws.Range("K110").Value = 2282.2856
$ws.Range("M110").Value = -237.2856000000002
$ws.Range("H122").Value = 3308.9707
$ws.Range("I122").Value = 3389.5715
$ws.Range("J122").Value = 2932.8333
$ws.Range("K122").Value = 10168.7145
$ws.Range("L122").Value = 8798.499899999999
$ws.Range("M122").Value = -7718.7145
$ws.Range("N122").Value = -13698.4999
$ws.Range("H132").Value = 3342.077
$ws.Range("I132").Value = 3342.077
$ws.Range("K132").Value = 10026.231
$ws.Range("M132").Value = -7496.231
$ws.Range("H136").Value = 2291.2903
$ws.Range("I136").Value = 2052.6155
$ws.Range("K136").Value = 6157.8465
$ws.Range("M136").Value = -3607.8465

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 656.7
$ws.Range("J80").Value = 1966.3334
$ws.Range("L80").Value = 1966.3334
$ws.Range("N80").Value = -3962.3334
$ws.Range("H83").Value = 656.7
$ws.Range("J83").Value = 1966.3334
$ws.Range("L83").Value = 9831.666999999999
$ws.Range("N83").Value = -19815.667
$ws.Range("H86").Value = 1526.5
$ws.Range("I86").Value = 1526.5
$ws.Range("K86").Value = 1526.5
$ws.Range("M86").Value = -403.5
$ws.Range("H89").Value = 1526.5
$ws.Range("I89").Value = 1526.5
$ws.Range("K89").Value = 7632.5
$ws.Range("M89").Value = -2016.5
$ws.Range("H105").Value = 1795.4584
$ws.Range("I105").Value = 1583.0555
$ws.Range("K105").Value = 1583.0555
$ws.Range("M105").Value = 163.9445000000001
$ws.Range("H107").Value = 1290.5
$ws.Range("I107").Value = 1250.8462
$ws.Range("K107").Value = 1250.8462
$ws.Range("M107").Value = 669.1538
$ws.Range("H134").Value = 2487.2144
$ws.Range("I134").Value = 2479.625
$ws.Range("J134").Value = 2497.3333
$ws.Range("K134").Value = 7438.875
$ws.Range("L134").Value = 7491.999899999999
$ws.Range("M134").Value = -4903.875
$ws.Range("N134").Value = -12561.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 556.4545000000001
$ws.Range("I22").Value = 542.1
$ws.Range("K22").Value = 542.1
$ws.Range("M22").Value = -192.1
$ws.Range("H31").Value = 4075.5217
$ws.Range("I31").Value = 2526.261
$ws.Range("J31").Value = 5624.7827
$ws.Range("K31").Value = 2526.261
$ws.Range("L31").Value = 5624.7827
$ws.Range("M31").Value = -2231.261
$ws.Range("N31").Value = -6214.7827
$ws.Range("H34").Value = 4075.5217
$ws.Range("I34").Value = 2526.261
$ws.Range("J34").Value = 5624.7827
$ws.Range("K34").Value = 2526.261
$ws.Range("L34").Value = 5624.7827
$ws.Range("M34").Value = -2324.261
$ws.Range("N34").Value = -6028.7827
$ws.Range("H58").Value = 3652.1904
$ws.Range("I58").Value = 3510.3684
$ws.Range("J58").Value = 4999.5
$ws.Range("K58").Value = 3510.3684
$ws.Range("L58").Value = 4999.5
$ws.Range("M58").Value = -3307.3684
$ws.Range("N58").Value = -5405.5
$ws.Range("H62").Value = 13600.083
$ws.Range("J62").Value = 8696
$ws.Range("L62").Value = 8696
$ws.Range("N62").Value = -9944
$ws.Range("H65").Value = 13600.083
$ws.Range("J65").Value = 8696
$ws.Range("L65").Value = 43480
$ws.Range("N65").Value = -49720
$ws.Range("H99").Value = 2339
$ws.Range("I99").Value = 2339
$ws.Range("K99").Value = 2339
$ws.Range("M99").Value = -841
$ws.Range("H107").Value = 515.1539
$ws.Range("I107").Value = 508.91666
$ws.Range("J107").Value = 590
$ws.Range("K107").Value = 508.91666
$ws.Range("L107").Value = 590
$ws.Range("M107").Value = 1411.08334
$ws.Range("N107").Value = -4430
$ws.Range("H122").Value = 4618.125
$ws.Range("J122").Value = 5436.6665
$ws.Range("L122").Value = 16309.9995
$ws.Range("N122").Value = -21209.9995
$ws.Range("H126").Value = 2339
$ws.Range("I126").Value = 2339
$ws.Range("K126").Value = 7017
$ws.Range("M126").Value = -4547
$ws.Range("H134").Value = 2550.9092
$ws.Range("I134").Value = 2577.1428
$ws.Range("K134").Value = 7731.428400000001
$ws.Range("M134").Value = -5196.428400000001
$ws.Range("H135").Value = 93852.664
$ws.Range("J135").Value = 93852.664
$ws.Range("L135").Value = 93852.664
$ws.Range("N135").Value = -103992.664
$ws.Range("H136").Value = 3652.1904
$ws.Range("I136").Value = 3510.3684
$ws.Range("J136").Value = 4999.5
$ws.Range("K136").Value = 10531.1052
$ws.Range("L136").Value = 14998.5
$ws.Range("M136").Value = -7981.1052
$ws.Range("N136").Value = -20098.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 274.25
$ws.Range("I5").Value = 282.33334
$ws.Range("J5").Value = 250
$ws.Range("K5").Value = 847.0000200000001
$ws.Range("L5").Value = 750
$ws.Range("M5").Value = -735.0000200000001
$ws.Range("N5").Value = -974
$ws.Range("H11").Value = 598.3333
$ws.Range("I11").Value = 319.5
$ws.Range("J11").Value = 737.75
$ws.Range("K11").Value = 958.5
$ws.Range("L11").Value = 2213.25
$ws.Range("M11").Value = -818.5
$ws.Range("N11").Value = -2493.25
$ws.Range("H37").Value = 81164.836
$ws.Range("J37").Value = 81164.836
$ws.Range("L37").Value = 243494.508
$ws.Range("N37").Value = -243718.508
$ws.Range("H64").Value = 1375
$ws.Range("I64").Value = 1100
$ws.Range("J64").Value = 1443.75
$ws.Range("K64").Value = 3300
$ws.Range("L64").Value = 4331.25
$ws.Range("M64").Value = -3030
$ws.Range("N64").Value = -4871.25
$ws.Range("H67").Value = 1375
$ws.Range("I67").Value = 1100
$ws.Range("J67").Value = 1443.75
$ws.Range("K67").Value = 3300
$ws.Range("L67").Value = 4331.25
$ws.Range("M67").Value = -2364
$ws.Range("N67").Value = -6203.25
$ws.Range("H75").Value = 475
$ws.Range("J75").Value = 450
$ws.Range("L75").Value = 1350
$ws.Range("N75").Value = -3346
$ws.Range("H78").Value = 475
$ws.Range("J78").Value = 450
$ws.Range("L78").Value = 4050
$ws.Range("N78").Value = -14034
$ws.Range("H87").Value = 25965.375
$ws.Range("I87").Value = 25965.375
$ws.Range("K87").Value = 77896.125
$ws.Range("M87").Value = -76648.125
$ws.Range("H90").Value = 25965.375
$ws.Range("I90").Value = 25965.375
$ws.Range("K90").Value = 233688.375
$ws.Range("M90").Value = -227448.375
$ws.Range("H94").Value = 1219.4286
$ws.Range("J94").Value = 1698.6
$ws.Range("L94").Value = 5095.799999999999
$ws.Range("N94").Value = -6447.799999999999
$ws.Range("H131").Value = 1142649.6
$ws.Range("I131").Value = 112027.6
$ws.Range("J131").Value = 1303684.4
$ws.Range("K131").Value = 336082.8
$ws.Range("L131").Value = 3911053.2
$ws.Range("M131").Value = -331042.8
$ws.Range("N131").Value = -3921133.2
$ws.Range("H135").Value = 274.25
$ws.Range("I135").Value = 282.33334
$ws.Range("J135").Value = 250
$ws.Range("K135").Value = 2541.00006
$ws.Range("L135").Value = 2250
$ws.Range("M135").Value = -6.000060000000303
$ws.Range("N135").Value = -7320

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8775.789000000001
$ws.Range("I70").Value = 8009.5
$ws.Range("J70").Value = 9627.223
$ws.Range("K70").Value = 8009.5
$ws.Range("L70").Value = 9627.223
$ws.Range("M70").Value = -7739.5
$ws.Range("N70").Value = -10167.223
$ws.Range("H73").Value = 8775.789000000001
$ws.Range("I73").Value = 8009.5
$ws.Range("J73").Value = 9627.223
$ws.Range("K73").Value = 8009.5
$ws.Range("L73").Value = 9627.223
$ws.Range("M73").Value = -7073.5
$ws.Range("N73").Value = -11499.223
$ws.Range("H80").Value = 7947.826
$ws.Range("I80").Value = 3668.5
$ws.Range("K80").Value = 3668.5
$ws.Range("M80").Value = -2670.5
$ws.Range("H83").Value = 7947.826
$ws.Range("I83").Value = 3668.5
$ws.Range("K83").Value = 18342.5
$ws.Range("M83").Value = -13350.5
$ws.Range("H97").Value = 636.2727
$ws.Range("I97").Value = 556
$ws.Range("J97").Value = 997.5
$ws.Range("K97").Value = 556
$ws.Range("L97").Value = 997.5
$ws.Range("M97").Value = -60
$ws.Range("N97").Value = -1989.5
$ws.Range("H102").Value = 2714.2222
$ws.Range("I102").Value = 2303.5
$ws.Range("K102").Value = 2303.5
$ws.Range("M102").Value = -681.5
$ws.Range("H107").Value = 442.66666
$ws.Range("I107").Value = 279.83334
$ws.Range("J107").Value = 489.1905
$ws.Range("K107").Value = 279.83334
$ws.Range("L107").Value = 489.1905
$ws.Range("M107").Value = 1640.16666
$ws.Range("N107").Value = -4329.1905
$ws.Range("H113").Value = 1887.421
$ws.Range("I113").Value = 1851.7858
$ws.Range("J113").Value = 1987.2
$ws.Range("K113").Value = 1851.7858
$ws.Range("L113").Value = 1987.2
$ws.Range("M113").Value = 318.2141999999999
$ws.Range("N113").Value = -6327.2
$ws.Range("H122").Value = 2382.8076
$ws.Range("I122").Value = 2302.8096
$ws.Range("J122").Value = 2718.8
$ws.Range("K122").Value = 6908.4288
$ws.Range("L122").Value = 8156.400000000001
$ws.Range("M122").Value = -4458.4288
$ws.Range("N122").Value = -13056.4
$ws.Range("H126").Value = 3297.9
$ws.Range("I126").Value = 3220.111
$ws.Range("K126").Value = 9660.332999999999
$ws.Range("M126").Value = -7190.332999999999
$ws.Range("I132").Value = 2063.818
$ws.Range("J132").Value = 2826.3333
$ws.Range("K132").Value = 6191.454000000001
$ws.Range("L132").Value = 8478.999899999999
$ws.Range("M132").Value = -3661.454000000001
$ws.Range("N132").Value = -13538.9999
$ws.Range("H133").Value = 76203.336
$ws.Range("J133").Value = 76203.336
$ws.Range("L133").Value = 76203.336
$ws.Range("N133").Value = -86323.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4299.5
$ws.Range("I7").Value = 4399.5
$ws.Range("J7").Value = 4249.5
$ws.Range("K7").Value = 4399.5
$ws.Range("L7").Value = 4249.5
$ws.Range("M7").Value = -4287.5
$ws.Range("N7").Value = -4473.5
$ws.Range("H16").Value = 576.46155
$ws.Range("I16").Value = 581.3333
$ws.Range("K16").Value = 581.3333
$ws.Range("M16").Value = -411.3333
$ws.Range("H22").Value = 1464.875
$ws.Range("I22").Value = 1060
$ws.Range("J22").Value = 1869.75
$ws.Range("K22").Value = 1060
$ws.Range("L22").Value = 1869.75
$ws.Range("M22").Value = -765
$ws.Range("N22").Value = -2459.75
$ws.Range("H27").Value = 1464.875
$ws.Range("I27").Value = 1060
$ws.Range("J27").Value = 1869.75
$ws.Range("K27").Value = 1060
$ws.Range("L27").Value = 1869.75
$ws.Range("M27").Value = -953
$ws.Range("N27").Value = -2083.75
$ws.Range("H68").Value = 965
$ws.Range("J68").Value = 965
$ws.Range("L68").Value = 965
$ws.Range("N68").Value = -2463
$ws.Range("H71").Value = 965
$ws.Range("J71").Value = 965
$ws.Range("L71").Value = 4825
$ws.Range("N71").Value = -12313
$ws.Range("H100").Value = 2812.5
$ws.Range("I100").Value = 2533.3333
$ws.Range("K100").Value = 2533.3333
$ws.Range("M100").Value = -1992.3333
$ws.Range("H122").Value = 5712.7617
$ws.Range("I122").Value = 5475
$ws.Range("J122").Value = 6473.6
$ws.Range("K122").Value = 16425
$ws.Range("L122").Value = 19420.8
$ws.Range("M122").Value = -13975
$ws.Range("N122").Value = -24320.8
$ws.Range("H126").Value = 4299.5
$ws.Range("I126").Value = 4399.5
$ws.Range("J126").Value = 4249.5
$ws.Range("K126").Value = 13198.5
$ws.Range("L126").Value = 12748.5
$ws.Range("M126").Value = -10728.5
$ws.Range("N126").Value = -17688.5
$ws.Range("H132").Value = 2947
$ws.Range("I132").Value = 2849.5293
$ws.Range("K132").Value = 8548.5879
$ws.Range("M132").Value = -6018.5879
$ws.Range("H136").Value = 3173.5334
$ws.Range("I136").Value = 2911.1667
$ws.Range("J136").Value = 4223
$ws.Range("K136").Value = 8733.500100000001
$ws.Range("L136").Value = 12669
$ws.Range("M136").Value = -6183.500100000001
$ws.Range("N136").Value = -17769

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 478.7
$ws.Range("I100").Value = 541.7143
$ws.Range("J100").Value = 331.66666
$ws.Range("K100").Value = 1083.4286
$ws.Range("L100").Value = 663.33332
$ws.Range("M100").Value = -542.4286
$ws.Range("N100").Value = -1745.33332
$ws.Range("H107").Value = 1639.0834
$ws.Range("I107").Value = 1715.5454
$ws.Range("J107").Value = 798
$ws.Range("K107").Value = 5146.6362
$ws.Range("L107").Value = 2394
$ws.Range("M107").Value = -3226.6362
$ws.Range("N107").Value = -6234
$ws.Range("H126").Value = 2940.6843
$ws.Range("I126").Value = 3058.6
$ws.Range("J126").Value = 2498.5
$ws.Range("K126").Value = 9175.799999999999
$ws.Range("L126").Value = 7495.5
$ws.Range("M126").Value = -6705.799999999999
$ws.Range("N126").Value = -12435.5
$ws.Range("H132").Value = 2102.1458
$ws.Range("I132").Value = 2060.2195
$ws.Range("J132").Value = 2347.7144
$ws.Range("K132").Value = 6180.6585
$ws.Range("L132").Value = 7043.1432
$ws.Range("M132").Value = -3650.6585
$ws.Range("N132").Value = -12103.1432
